$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Session 4 (Greedy) mark for the student
$ws.Range("E4").Value = 2.5

# Feedback comment for Session 4 (Greedy)
$ws.Range("E5").Value = "Algorithms didn't give correct results and the expected complexity for greedy 2 and greedy 3 is O(nlogn) if you sort the elements beforehand or if you use a priority queue to have a greedy algorithm"

# Update the selection to reflect the newly edited column/range
$ws.Range("E5:E12").Select()
